# PaymentReceipt.docx: rework the "received from ... {lawyerName}" sentence
# to add a "by {paymentMethod}" clause after the (renamed) deposit-amount
# placeholder, and lower-case the {Date} -> {date} merge field.
#
# Before:
#   ... received from {clientName} the sum of $ {amount} {paymentMethod}
#   for {reason} with Me {lawyerName}.
#   Montreal, {Date}
#
# After:
#   ... received from {clientName} the sum of $ {depositAmount} by
#   {paymentMethod} for {reason} with Me {lawyerName}.
#   Montreal, {date}

$d = $word.ActiveDocument

# {clientName} was split across "{" / "clientName" / "}" runs (with
# spell-check proofErr markers); normalize it to a single clean run.
$d.Content.Find.Execute("{clientName}", $true, $false, $false, $false, $false, `
    $true, 1, $false, "{clientName}", 2) | Out-Null

# {amount} -> {depositAmount} (stays inside the underlined "$ {...}" run)
$d.Content.Find.Execute("{amount}", $true, $false, $false, $false, $false, `
    $true, 1, $false, "{depositAmount}", 2) | Out-Null

# Insert "by " before "{paymentMethod}" (plain, non-underlined text)
$d.Content.Find.Execute(" {paymentMethod} for {reason} with ", $true, $false, `
    $false, $false, $false, $true, 1, $false, `
    " by {paymentMethod} for {reason} with ", 2) | Out-Null

# "Me {lawyerName}" was split across several runs (with spell-check
# proofErr markers around "lawyerName"); normalize it to clean runs.
$d.Content.Find.Execute("Me {lawyerName}", $true, $false, $false, $false, `
    $false, $true, 1, $false, "Me {lawyerName}", 2) | Out-Null

# {Date} -> {date}
$d.Content.Find.Execute("{Date}", $true, $false, $false, $false, $false, `
    $true, 1, $false, "{date}", 2) | Out-Null
